$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Starting state:
#   Para 1: [bookmarkStart/_GoBack][bookmarkEnd] "Assignment"
#   Para 2: "CSA0805 \u2013 Python Programming"
#
# Target state:
#   Para 1: "Assignment"                                   (bookmark removed)
#   Para 2: "CSA08" + "14" (distinct run) + [_GoBack] + " \u2013 Python Programming"
# ---------------------------------------------------------------------------

# 1) Drop the old _GoBack bookmark from paragraph 1 (it gets re-inserted later,
#    right after the new "14" run in paragraph 2).
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# 2) Locate "05" inside "CSA0805" (the course-code suffix being updated to "14").
$target = $d.Content
$target.Find.ClearFormatting()
$target.Find.Execute("CSA0805", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$codeStart = $target.Start
$digitsRange = $d.Range($codeStart + 5, $codeStart + 7)
Write-Output ("digits range text=[" + $digitsRange.Text + "]")

# 3) Replace "05" with "14" in place.
$digitsRange.Text = "14"
$newDigits = $d.Range($codeStart + 5, $codeStart + 7)
Write-Output ("new digits text=[" + $newDigits.Text + "]")

# 4) Force the "14" text into its own run (distinct rPr) by nudging a
#    character-level property on just that sub-range.
$newDigits.Italic = 1
$newDigits.Italic = 0

# 5) Re-insert the _GoBack bookmark immediately after "14" (collapsed range).
$afterDigits = $codeStart + 7
$bmRange = $d.Range($afterDigits, $afterDigits)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

Write-Output ("Para2 text now=[" + $d.Paragraphs.Item(2).Range.Text + "]")
